$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Muskan Vaswan", "2021-01-15 07:08:54.649707"),
    @("Muskan Vaswan", "2021-01-15 10:19:42.224172"),
    @("Muskan Vaswan", "2021-01-15 10:20:36.332329"),
    @("Muskan Vaswan", "2021-01-15 10:21:41.973739"),
    @("Muskan Vaswan", "2021-01-15 10:22:13.154482"),
    @("Muskan Vaswan", "2021-01-15 10:22:56.788287"),
    @("Muskan Vaswan", "2021-01-15 10:23:10.748782")
)

$startRow = 32
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
